$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 updates: fill in the previously-missing trade-close data ---
# B3: Profitable = TRUE
$ws.Range("B3").Value = $true

# E3: SellPrice
$ws.Range("E3").Value = 313.26998900000001

# F3: Price Change %
$ws.Range("F3").Value = 2.1788019831044831

# G3: Holding flips from TRUE (still holding) to FALSE (position closed)
$ws.Range("G3").Value = $false

# --- Row 4 (new): next trade entry, only Principle recorded so far ---
$ws.Range("C4").Value = 10029.69

# --- Column width adjustments (bestFit recompute for the new, wider content) ---
$ws.Columns.Item(3).ColumnWidth = 8
$ws.Columns.Item(5).ColumnWidth = 10
